$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 13 ("Frontend Server" architecture diagram): add two red
#    "a)" / "b)" label text boxes next to the two message-flow arrows.
# ---------------------------------------------------------------------
$s13 = $p.Slides.Item(13)

$tbA = $s13.Shapes.AddTextbox(1, 549.4470866141733, 82.09614173228347, 58.84614173228346, 60.58590551181102)
$tbA.Name = "CasellaDiTesto 29"
$tbA.TextFrame.WordWrap = -1
$tbA.TextFrame.AutoSize = 1
$tbA.Fill.Visible = 0
$trA = $tbA.TextFrame.TextRange
$trA.Text = "a)"
$trA.Font.Size = 44
$trA.Font.Color.RGB = 255

$tbB = $s13.Shapes.AddTextbox(1, 509.0761417322835, 303.08551181102365, 58.84614173228346, 60.58590551181102)
$tbB.Name = "CasellaDiTesto 37"
$tbB.TextFrame.WordWrap = -1
$tbB.TextFrame.AutoSize = 1
$tbB.Fill.Visible = 0
$trB = $tbB.TextFrame.TextRange
$trB.Text = "b)"
$trB.Font.Size = 44
$trB.Font.Color.RGB = 255

# ---------------------------------------------------------------------
# 2) Slide 5 (state diagram): remove the dangling dashed connector that
#    used to leave "handleFront" toward nothing useful.
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item("Connettore 2 5").Delete()

# ---------------------------------------------------------------------
# 3) Slide 6 (state diagram continuation): highlight "handleFront" in
#    red (it is the entry point reached from slide 5) and clean up the
#    now-redundant curved connector / its "frontSonar" label plus the
#    stray dashed incoming connector.
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item("Ovale 4").Line.ForeColor.RGB = 255

$s6.Shapes.Item("Connettore curvo 268").Delete()
$s6.Shapes.Item("CasellaDiTesto 300").Delete()
$s6.Shapes.Item("Connettore 2 2").Delete()
